$d = $word.ActiveDocument
$euro = [char]0x20AC

# --- Edit 1: "within the RPA an amount ..." -> "within the RPA is an amount ..." ---
# Insert the word "is " right after "within the RPA ".
$d.Content.Find.Execute(
    "within the RPA an amount of k" + $euro + " 1,190 available",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "within the RPA is an amount of k" + $euro + " 1,190 available",
    2
)

# --- Edit 2: ", with k€ 60/20 additional for experimental/computational work" ---
#             -> ", with k€ 60 (or 20) additional for experimental (or computational) work"
# a) "60/20" -> "60 (or 20)"
$d.Content.Find.Execute(
    "k" + $euro + " 60/20 additional",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "k" + $euro + " 60 (or 20) additional",
    2
)

# b) "experimental/computational" -> "experimental (or computational)"
$d.Content.Find.Execute(
    "experimental/computational work",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "experimental (or computational) work",
    2
)
